$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds numeric-looking price text; prefix with an apostrophe so
# Excel stores it as text (matching the inline-string source data) instead
# of auto-converting it to a number.

$ws.Range("D2").Value = "'27.856.15"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "'1.749.54"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'332.98"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("D7").Value = "'0.3870"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "'0.3387"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'45.78"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").Value = "'1.115"
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").Value = "'0.07201"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "'22.48"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'6.177"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("D15").Value = "'1.748.96"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").Value = "'7.076"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'0.00001059"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "'0.06604"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'79.29"
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "'16.77"
$ws.Range("E21").Value = "  -3.19%  "
$ws.Range("D22").Value = "'6.183"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "'27.870.42"
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'11.65"
$ws.Range("E24").Value = "  -3.22%  "
$ws.Range("D25").Value = "'2.407"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("D26").Value = "'154.04"
$ws.Range("E26").Value = "  +1.10%  "
$ws.Range("D27").Value = "'19.84"
$ws.Range("E27").Value = "  -3.82%  "
$ws.Range("D28").Value = "'2.298"
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").Value = "'1.950.02"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'1.280"
$ws.Range("E30").Value = "  -10.37%  "
$ws.Range("D31").Value = "'130.92"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").Value = "'4.027"
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("D33").Value = "'5.822"
$ws.Range("E33").Value = "  -4.69%  "
$ws.Range("D34").Value = "'0.08792"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'12.12"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").Value = "'1.536"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "'0.6527"
$ws.Range("E37").Value = "  -3.81%  "

# Rows 38/39: VeChain and InternetComputer(DFINITY) swap positions, each with updated Price/Volume
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.131"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02276"
$ws.Range("E39").Value = "  -5.80%  "

$ws.Range("D40").Value = "'0.06117"
$ws.Range("E40").Value = "  -3.01%  "
$ws.Range("D41").Value = "'0.2100"
$ws.Range("E41").Value = "  -3.60%  "
$ws.Range("D42").Value = "'1.205"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").Value = "'8.011"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'13.74"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "'3.812"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'0.6035"
$ws.Range("E47").Value = "  -3.50%  "
$ws.Range("D48").Value = "'126.96"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").Value = "'1.991"
$ws.Range("E49").Value = "  -3.87%  "
$ws.Range("D50").Value = "'1.161"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("E51").Value = "  +3.65%  "
